$d = $word.ActiveDocument

$replacements = @(
    @("735÷4=183, 3", "643÷5=128, 3"),
    @("675÷4=168, 3", "914÷8=114, 2"),
    @("955÷5=191, 0", "646÷5=129, 1"),
    @("963÷5=192, 3", "306÷8=38, 2"),
    @("354÷9=39, 3", "502÷8=62, 6"),
    @("722÷7=103, 1", "882÷8=110, 2"),
    @("615÷9=68, 3", "405÷9=45, 0"),
    @("352÷7=50, 2", "416÷5=83, 1"),
    @("643÷7=91, 6", "464÷7=66, 2"),
    @("855÷6=142, 3", "906÷9=100, 6"),
    @("256÷4=64, 0", "145÷9=16, 1"),
    @("914÷4=228, 2", "751÷3=250, 1"),
    @("341÷9=37, 8", "104÷2=52, 0"),
    @("504÷4=126, 0", "981÷7=140, 1"),
    @("133÷3=44, 1", "343÷7=49, 0"),
    @("639÷5=127, 4", "591÷7=84, 3"),
    @("523÷8=65, 3", "655÷7=93, 4"),
    @("731÷4=182, 3", "417÷8=52, 1"),
    @("368÷8=46, 0", "539÷4=134, 3"),
    @("211÷6=35, 1", "460÷5=92, 0"),
    @("323÷3=107, 2", "641÷3=213, 2"),
    @("329÷2=164, 1", "190÷2=95, 0"),
    @("435÷9=48, 3", "495÷5=99, 0"),
    @("940÷7=134, 2", "258÷5=51, 3"),
    @("284÷3=94, 2", "782÷5=156, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying replacements"
